$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.715.36'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.847.19'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.69'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4308'
$ws.Range('E7').Value = '  +1.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3653'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07346'
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8762'
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').Value = '1.904.07'
$ws.Range('E12').Value = '  +4.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.353'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.529'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06959'
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.80'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009003'
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.37'
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('D21').Value = '27.796.83'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.983'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.34'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').Value = '2.062.46'
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.981'
$ws.Range('E25').Value = '  -3.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.12'
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.62'
$ws.Range('E27').Value = '  +2.11%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '119.99'
$ws.Range('E28').Value = '  +8.06%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.259'
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.875'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08902'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7563'
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.555'
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.958'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.128'
$ws.Range('E35').Value = '  +2.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.108'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05431'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01934'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.824'
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5092'
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1662'
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.582'
$ws.Range('E43').Value = '  -3.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.386'
$ws.Range('E44').Value = '  +1.88%  '
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.06543'
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.34'
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4668'
$ws.Range('E48').Value = '  -0.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.0000'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.631'
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.74'
$ws.Range('E51').Value = '  +0.47%  '
